# Update "想去人数" (want-to-go count) figures for three events.
# The workbook has the same events duplicated across the "展览" sheet
# (category-specific) and the "全部类型" sheet (all categories combined),
# so both copies need to be updated to keep the data consistent.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F3").Value = 2068
        $ws.Range("F5").Value = 1190
        $ws.Range("F6").Value = 354
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F3").Value = 2068
        $ws.Range("F7").Value = 1190
        $ws.Range("F8").Value = 354
    }
}
